$wb = $excel.ActiveWorkbook

# --- Sheet "SoCDTtiNTY-psgr": LDVs row (row 2) values 0.0755 -> 0.076 ---
$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsPsgr.Range("B2:H2").Value = 0.076
$wsPsgr.Rows.Item(1).RowHeight = 30

# --- Sheet "SoCDTtiNTY-frgt": HDVs row (row 3) values -> uniform 0.035 ---
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")
$wsFrgt.Range("B3:H3").Value = 0.035
$wsFrgt.Rows.Item(1).RowHeight = 30

# The "About" sheet ends up as the active/selected tab in the saved file
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
